$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Insert a brand-new worksheet for "2022-Q4" right after "总计",
#    copying the layout/formatting of the "2022-Q3" sheet so the new
#    tab looks like its siblings, then overwrite its data.
# ------------------------------------------------------------------
$summarySheet = $wb.Worksheets.Item("总计")
$templateSheet = $wb.Worksheets.Item("2022-Q3")
$templateSheet.Copy($null, $summarySheet)

$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# fund code / name (B2,C2) stay the same as the template; update the
# rest of the metrics for the new quarter.
$q4.Cells.Item(2,4).Value = "'0.64"
$q4.Cells.Item(2,5).Value = "'93.56"
$q4.Cells.Item(2,6).Value = "'5.54"
$q4.Cells.Item(2,7).Value = "'0.0355"
$q4.Cells.Item(2,8).Value = 4

# ------------------------------------------------------------------
# 2) Insert a new row 2 into "总计" for the "2022-Q4" summary entry.
#    Everything below shifts down one row automatically.
# ------------------------------------------------------------------
$summarySheet.Rows.Item(2).Insert()
$summarySheet.Cells.Item(2,1).Value = 0
$summarySheet.Cells.Item(2,2).Value = "2022-Q4"
$summarySheet.Cells.Item(2,3).Value = 1
$summarySheet.Cells.Item(2,4).Value = 0.04

# Match formatting of the rest of the data rows for the new row.
$summarySheet.Cells.Item(2,1).Font.Bold = $true
$summarySheet.Cells.Item(2,1).Borders.LineStyle = 1

# The row that fell off the bottom (old row 9, "2020-Q4") is now row
# 10 thanks to the insert above; its index column wasn't auto-updated
# by the shift, so fix it up explicitly.
$summarySheet.Cells.Item(10,1).Value = 8

# ------------------------------------------------------------------
# 3) Restore "2020-Q4" as the selected/active tab (it was the active
#    tab before the edit; copying a sheet above made the copy active).
# ------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
